$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the new team-record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, bordered, centered) from an existing header cell
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the team record values for every data row (2 through 48)
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 66
    $ws.Cells.Item($row, 31).Value = 96
    $ws.Cells.Item($row, 32).Value = 0
}
